# Scheduled-runner market data refresh: pushes newly-fetched Market Board
# averages/prices (and recomputed profit deltas) into a handful of Leve
# rows across each crafting-class worksheet. Only numeric columns
# H (currentAveragePrice), I (currentAveragePriceNQ), J (currentAveragePriceHQ),
# K (LevePriceNQ), L (LevePriceHQ), M (LeveProfitNQ) and N (LeveProfitHQ)
# are touched; a few rows gain a previously-absent M cell once NQ profit
# becomes computable again.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 61624.65
$ws.Range("J64").Value = 3050
$ws.Range("L64").Value = 3050
$ws.Range("N64").Value = -3546
$ws.Range("H67").Value = 61624.65
$ws.Range("J67").Value = 3050
$ws.Range("L67").Value = 3050
$ws.Range("N67").Value = -4766
$ws.Range("H112").Value = 1282.15
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1302.2051
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 3906.615299999999
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -6122.615299999999
$ws.Range("H129").Value = 2427.8262
$ws.Range("J129").Value = 1056.8448
$ws.Range("L129").Value = 3170.5344
$ws.Range("N129").Value = -13170.5344
$ws.Range("H137").Value = 1796.6364
$ws.Range("I137").Value = 1523.4348
$ws.Range("J137").Value = 2425
$ws.Range("K137").Value = 4570.3044
$ws.Range("L137").Value = 7275
$ws.Range("M137").Value = -2020.3044
$ws.Range("N137").Value = -12375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 471.85715
$ws.Range("I5").Value = 567
$ws.Range("J5").Value = 400.5
$ws.Range("K5").Value = 567
$ws.Range("L5").Value = 400.5
$ws.Range("M5").Value = -455
$ws.Range("N5").Value = -624.5
$ws.Range("H32").Value = 28944.902
$ws.Range("I32").Value = 8287.954
$ws.Range("J32").Value = 256171.33
$ws.Range("K32").Value = 8287.954
$ws.Range("L32").Value = 256171.33
$ws.Range("M32").Value = -8000.954
$ws.Range("N32").Value = -256745.33
$ws.Range("H63").Value = 2800
$ws.Range("I63").Value = 2200
$ws.Range("J63").Value = 3400
$ws.Range("K63").Value = 2200
$ws.Range("L63").Value = 3400
$ws.Range("M63").Value = -1514
$ws.Range("N63").Value = -4772
$ws.Range("H66").Value = 2800
$ws.Range("I66").Value = 2200
$ws.Range("J66").Value = 3400
$ws.Range("K66").Value = 11000
$ws.Range("L66").Value = 17000
$ws.Range("M66").Value = -7568
$ws.Range("N66").Value = -23864
$ws.Range("H74").Value = 896.9091
$ws.Range("I74").Value = 795.2941
$ws.Range("K74").Value = 795.2941
$ws.Range("M74").Value = 78.70590000000004
$ws.Range("H77").Value = 896.9091
$ws.Range("I77").Value = 795.2941
$ws.Range("K77").Value = 3976.4705
$ws.Range("M77").Value = 391.5295000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 471.85715
$ws.Range("I4").Value = 567
$ws.Range("J4").Value = 400.5
$ws.Range("K4").Value = 567
$ws.Range("L4").Value = 400.5
$ws.Range("M4").Value = -452
$ws.Range("N4").Value = -630.5
$ws.Range("H35").Value = 19684.666
$ws.Range("J35").Value = 19684.666
$ws.Range("L35").Value = 19684.666
$ws.Range("N35").Value = -20304.666
$ws.Range("H82").Value = 21499
$ws.Range("I82").Value = 1752.3334
$ws.Range("J82").Value = 29961.857
$ws.Range("K82").Value = 1752.3334
$ws.Range("L82").Value = 29961.857
$ws.Range("M82").Value = -1369.3334
$ws.Range("N82").Value = -30727.857
$ws.Range("H85").Value = 21499
$ws.Range("I85").Value = 1752.3334
$ws.Range("J85").Value = 29961.857
$ws.Range("K85").Value = 1752.3334
$ws.Range("L85").Value = 29961.857
$ws.Range("M85").Value = -426.3334
$ws.Range("N85").Value = -32613.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 18000
$ws.Range("J95").Value = 18000
$ws.Range("L95").Value = 18000
$ws.Range("N95").Value = -23492
$ws.Range("H132").Value = 3411.3684
$ws.Range("I132").Value = 3301.2
$ws.Range("J132").Value = 3824.5
$ws.Range("K132").Value = 9903.599999999999
$ws.Range("L132").Value = 11473.5
$ws.Range("M132").Value = -7373.599999999999
$ws.Range("N132").Value = -16533.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1361.8
$ws.Range("I5").Value = 1523
$ws.Range("J5").Value = 1325.9778
$ws.Range("K5").Value = 4569
$ws.Range("L5").Value = 3977.9334
$ws.Range("M5").Value = -4457
$ws.Range("N5").Value = -4201.9334
$ws.Range("H58").Value = 2000
$ws.Range("J58").Value = 1733.3334
$ws.Range("L58").Value = 5200.0002
$ws.Range("N58").Value = -5456.0002
$ws.Range("H68").Value = 457.14285
$ws.Range("I68").Value = 543.3333
$ws.Range("J68").Value = 392.5
$ws.Range("K68").Value = 1629.9999
$ws.Range("L68").Value = 1177.5
$ws.Range("M68").Value = -818.9999
$ws.Range("N68").Value = -2799.5
$ws.Range("H71").Value = 457.14285
$ws.Range("I71").Value = 543.3333
$ws.Range("J71").Value = 392.5
$ws.Range("K71").Value = 4889.9997
$ws.Range("L71").Value = 3532.5
$ws.Range("M71").Value = -833.9997000000003
$ws.Range("N71").Value = -11644.5
$ws.Range("H101").Value = 4119.6
$ws.Range("J101").Value = 4119.6
$ws.Range("L101").Value = 12358.8
$ws.Range("N101").Value = -17226.8
$ws.Range("H131").Value = 774626.4399999999
$ws.Range("I131").Value = 525.75
$ws.Range("J131").Value = 946648.8
$ws.Range("K131").Value = 1577.25
$ws.Range("L131").Value = 2839946.4
$ws.Range("M131").Value = 3462.75
$ws.Range("N131").Value = -2850026.4
$ws.Range("H135").Value = 1361.8
$ws.Range("I135").Value = 1523
$ws.Range("J135").Value = 1325.9778
$ws.Range("K135").Value = 13707
$ws.Range("L135").Value = 11933.8002
$ws.Range("M135").Value = -11172
$ws.Range("N135").Value = -17003.8002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 2250.2727
$ws.Range("I17").Value = 218
$ws.Range("J17").Value = 7669.6665
$ws.Range("K17").Value = 218
$ws.Range("L17").Value = 7669.6665
$ws.Range("M17").Value = -50
$ws.Range("N17").Value = -8005.6665
$ws.Range("H80").Value = 91110900
$ws.Range("I80").Value = 167034830
$ws.Range("J80").Value = 2170
$ws.Range("K80").Value = 167034830
$ws.Range("L80").Value = 2170
$ws.Range("M80").Value = -167033832
$ws.Range("N80").Value = -4166
$ws.Range("H83").Value = 91110900
$ws.Range("I83").Value = 167034830
$ws.Range("J83").Value = 2170
$ws.Range("K83").Value = 835174150
$ws.Range("L83").Value = 10850
$ws.Range("M83").Value = -835169158
$ws.Range("N83").Value = -20834
$ws.Range("H119").Value = 39993.332
$ws.Range("J119").Value = 39993.332
$ws.Range("L119").Value = 39993.332
$ws.Range("N119").Value = -49669.332
$ws.Range("H132").Value = 3845.8845
$ws.Range("I132").Value = 3070.1428
$ws.Range("J132").Value = 4750.9165
$ws.Range("K132").Value = 9210.428400000001
$ws.Range("L132").Value = 14252.7495
$ws.Range("M132").Value = -6680.428400000001
$ws.Range("N132").Value = -19312.7495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 75180.25
$ws.Range("J13").Value = 75180.25
$ws.Range("L13").Value = 75180.25
$ws.Range("N13").Value = -75460.25
$ws.Range("H132").Value = 3069.9312
$ws.Range("I132").Value = 3301.077
$ws.Range("J132").Value = 1066.6666
$ws.Range("K132").Value = 9903.231
$ws.Range("L132").Value = 3199.9998
$ws.Range("M132").Value = -7373.231
$ws.Range("N132").Value = -8259.9998
$ws.Range("H133").Value = 34544
$ws.Range("J133").Value = 34544
$ws.Range("L133").Value = 34544
$ws.Range("N133").Value = -39604

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 39830
$ws.Range("J117").Value = 39830
$ws.Range("L117").Value = 39830
$ws.Range("N117").Value = -49008
$ws.Range("H132").Value = 6775.433
$ws.Range("I132").Value = 3964.842
$ws.Range("J132").Value = 11630.091
$ws.Range("K132").Value = 11894.526
$ws.Range("L132").Value = 34890.273
$ws.Range("M132").Value = -9364.526
$ws.Range("N132").Value = -39950.273
